# Update the multiplication-fact table with the newly generated answers.
# The single table in the document holds 5 data rows (1, 5, 10, 15, 20 in
# Word's 1-based row numbering), each with 5 columns. We address each cell
# directly via Table.Cell(row, col) rather than doing a global text
# find/replace, because some new values coincide with other cells' old
# values (e.g. "750x6=4500" is both an old value in row 20 and the new
# value for row 10), which would make a sequential text-replace ambiguous.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("626×9=5634", "485×6=2910", "304×5=1520", "656×8=5248", "266×3=798")
    5  = @("782×5=3910", "177×6=1062", "376×7=2632", "715×2=1430", "613×7=4291")
    10 = @("750×6=4500", "542×3=1626", "126×2=252",  "537×5=2685", "559×6=3354")
    15 = @("832×6=4992", "500×4=2000", "147×7=1029", "952×8=7616", "457×3=1371")
    20 = @("904×9=8136", "114×4=456",  "747×3=2241", "994×5=4970", "172×4=688")
}

foreach ($rowIndex in $newValues.Keys) {
    $rowValues = $newValues[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $rowValues[$col - 1]
    }
}
